$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 13.355
$ws.Range("E6").Value = 13.228
$ws.Range("C7").Value = -13.498
$ws.Range("A10").Value = -20.926
$ws.Range("A12").Value = -21.808
$ws.Range("B13").Value = 6.724000000000001
$ws.Range("A18").Value = -21.808
$ws.Range("C20").Value = -13.041
